$wb = $excel.ActiveWorkbook

# 1. Update version number on the "isa_template" sheet
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.4"

# 2. Update building block header names and URL values on "plant_source" sheet
$wsPlant = $wb.Worksheets.Item("plant_source")

$wsPlant.Range("Q1").Value = "Characteristic [geographic location]"
$wsPlant.Range("T1").Value = "Characteristic [latitude]"
$wsPlant.Range("W1").Value = "Characteristic [longitude]"

# 3. Update ontology URL/term values in the data row
$wsPlant.Range("D2").Value = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_1000169"
$wsPlant.Range("S2").Value = "https://bioregistry.io/NCIT:C16636"
$wsPlant.Range("AB2").Value = "https://www.ebi.ac.uk/ols4/ontologies/po/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FPO_0025034"
$wsPlant.Range("AH2").Value = "https://bioregistry.io/EFO:0005168"
